# The deck ships two theme parts:
#   - theme2.xml ("Integral" / "Red Violet") is the theme actually bound to
#     the slide master (and therefore to every slide) via Presentation.Designs(1).
#   - theme1.xml ("Office Theme" / "Office") is only referenced by the notes
#     master and isn't reachable as an independent, editable theme through the
#     PowerPoint object model (NotesMaster/HandoutMaster/TitleMaster all expose
#     the same single live Theme as the slide master).
#
# The authored change swaps the two themes' contents (theme1 <-> theme2).
# The only real content difference between them is the 12 color-scheme
# entries (font scheme / format scheme are identical byte-for-byte), so the
# faithful, COM-reachable equivalent is to repaint the live theme's color
# scheme from "Red Violet" over to the "Office" palette that the swap puts
# into the slide master's theme file.

$p = $ppt.ActivePresentation
$design = $p.Designs.Item(1)
$master = $design.SlideMaster
$colorScheme = $master.Theme.ThemeColorScheme

# Office theme palette, expressed as the OLE RGB() integer (R + G*256 + B*65536):
#  1 dk1       000000 ->        0
#  2 lt1       FFFFFF -> 16777215
#  3 dk2       44546A ->  6968388
#  4 lt2       E7E6E6 -> 15132391
#  5 accent1   5B9BD5 -> 13998939
#  6 accent2   ED7D31 ->  3243501
#  7 accent3   A5A5A5 -> 10855845
#  8 accent4   FFC000 ->    49407
#  9 accent5   4472C4 -> 12874308
# 10 accent6   70AD47 ->  4697456
# 11 hlink     0563C1 -> 12673797
# 12 folHlink  954F72 ->  7491477
$colorScheme.Item(1).RGB  = 0
$colorScheme.Item(2).RGB  = 16777215
$colorScheme.Item(3).RGB  = 6968388
$colorScheme.Item(4).RGB  = 15132391
$colorScheme.Item(5).RGB  = 13998939
$colorScheme.Item(6).RGB  = 3243501
$colorScheme.Item(7).RGB  = 10855845
$colorScheme.Item(8).RGB  = 49407
$colorScheme.Item(9).RGB  = 12874308
$colorScheme.Item(10).RGB = 4697456
$colorScheme.Item(11).RGB = 12673797
$colorScheme.Item(12).RGB = 7491477
